$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4571.857
$ws.Range("I74").Value = 4334.3335
$ws.Range("J74").Value = 4750
$ws.Range("K74").Value = 4334.3335
$ws.Range("L74").Value = 4750
$ws.Range("M74").Value = -3398.3335
$ws.Range("N74").Value = -6622
$ws.Range("H77").Value = 4571.857
$ws.Range("I77").Value = 4334.3335
$ws.Range("J77").Value = 4750
$ws.Range("K77").Value = 21671.6675
$ws.Range("L77").Value = 23750
$ws.Range("M77").Value = -16991.6675
$ws.Range("N77").Value = -33110
$ws.Range("H86").Value = 32535.154
$ws.Range("I86").Value = 905.8889
$ws.Range("J86").Value = 103701
$ws.Range("K86").Value = 905.8889
$ws.Range("L86").Value = 103701
$ws.Range("M86").Value = 217.1111
$ws.Range("N86").Value = -105947
$ws.Range("H89").Value = 32535.154
$ws.Range("I89").Value = 905.8889
$ws.Range("J89").Value = 103701
$ws.Range("K89").Value = 4529.444500000001
$ws.Range("L89").Value = 518505
$ws.Range("M89").Value = 1086.555499999999
$ws.Range("N89").Value = -529737
$ws.Range("H92").Value = 217.92
$ws.Range("I92").Value = 189
$ws.Range("K92").Value = 189
$ws.Range("M92").Value = 1059
$ws.Range("H113").Value = 7845.2085
$ws.Range("I113").Value = 2957.8333
$ws.Range("J113").Value = 12732.583
$ws.Range("K113").Value = 2957.8333
$ws.Range("L113").Value = 12732.583
$ws.Range("M113").Value = 296.1667000000002
$ws.Range("N113").Value = -19240.583
$ws.Range("H129").Value = 1467.5294
$ws.Range("J129").Value = 1607.0714
$ws.Range("L129").Value = 4821.2142
$ws.Range("N129").Value = -14821.2142
$ws.Range("H137").Value = 1706.8214
$ws.Range("I137").Value = 1394.909
$ws.Range("J137").Value = 2850.5
$ws.Range("K137").Value = 4184.727000000001
$ws.Range("L137").Value = 8551.5
$ws.Range("M137").Value = -1634.727000000001
$ws.Range("N137").Value = -13651.5
$ws.Range("H138").Value = 1858.3191
$ws.Range("I138").Value = 1021.1053
$ws.Range("J138").Value = 2426.4285
$ws.Range("K138").Value = 3063.3159
$ws.Range("L138").Value = 7279.2855
$ws.Range("M138").Value = 2076.6841
$ws.Range("N138").Value = -17559.2855

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 50002
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H23").Value = 31068.834
$ws.Range("I23").Value = 70006
$ws.Range("J23").Value = 23281.4
$ws.Range("K23").Value = 70006
$ws.Range("L23").Value = 23281.4
$ws.Range("M23").Value = -69747
$ws.Range("N23").Value = -23799.4
$ws.Range("H32").Value = 7715.486
$ws.Range("I32").Value = 4098.129
$ws.Range("J32").Value = 35750
$ws.Range("K32").Value = 4098.129
$ws.Range("L32").Value = 35750
$ws.Range("M32").Value = -3811.129
$ws.Range("N32").Value = -36324
$ws.Range("H61").Value = 9072.666999999999
$ws.Range("I61").Value = 10000
$ws.Range("K61").Value = 10000
$ws.Range("M61").Value = -9788
$ws.Range("H97").Value = 1390.4
$ws.Range("I97").Value = 1254.7826
$ws.Range("K97").Value = 1254.7826
$ws.Range("M97").Value = -758.7826
$ws.Range("H132").Value = 1916.7843
$ws.Range("I132").Value = 768.5294
$ws.Range("J132").Value = 4213.294
$ws.Range("K132").Value = 2305.5882
$ws.Range("L132").Value = 12639.882
$ws.Range("M132").Value = 224.4117999999999
$ws.Range("N132").Value = -17699.882
$ws.Range("H136").Value = 9072.666999999999
$ws.Range("I136").Value = 10000
$ws.Range("K136").Value = 30000
$ws.Range("M136").Value = -27450

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5243.1333
$ws.Range("I20").Value = 2071.75
$ws.Range("J20").Value = 8867.571
$ws.Range("K20").Value = 2071.75
$ws.Range("L20").Value = 8867.571
$ws.Range("M20").Value = -1824.75
$ws.Range("N20").Value = -9361.571
$ws.Range("H22").Value = 6356.5
$ws.Range("I22").Value = 7250.2856
$ws.Range("K22").Value = 7250.2856
$ws.Range("M22").Value = -7077.2856
$ws.Range("H134").Value = 2915.45
$ws.Range("I134").Value = 1620.7142
$ws.Range("J134").Value = 3612.6155
$ws.Range("K134").Value = 4862.142599999999
$ws.Range("L134").Value = 10837.8465
$ws.Range("M134").Value = -2327.142599999999
$ws.Range("N134").Value = -15907.8465

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 180.75
$ws.Range("I7").Value = 55.5
$ws.Range("J7").Value = 205.8
$ws.Range("K7").Value = 55.5
$ws.Range("L7").Value = 205.8
$ws.Range("M7").Value = 57.5
$ws.Range("N7").Value = -431.8
$ws.Range("H31").Value = 1768.1277
$ws.Range("I31").Value = 1371.1052
$ws.Range("J31").Value = 3444.4443
$ws.Range("K31").Value = 1371.1052
$ws.Range("L31").Value = 3444.4443
$ws.Range("M31").Value = -1076.1052
$ws.Range("N31").Value = -4034.4443
$ws.Range("H34").Value = 1768.1277
$ws.Range("I34").Value = 1371.1052
$ws.Range("J34").Value = 3444.4443
$ws.Range("K34").Value = 1371.1052
$ws.Range("L34").Value = 3444.4443
$ws.Range("M34").Value = -1169.1052
$ws.Range("N34").Value = -3848.4443
$ws.Range("H58").Value = 1404.5122
$ws.Range("I58").Value = 703.7727
$ws.Range("J58").Value = 2215.8948
$ws.Range("K58").Value = 703.7727
$ws.Range("L58").Value = 2215.8948
$ws.Range("M58").Value = -500.7727
$ws.Range("N58").Value = -2621.8948
$ws.Range("H132").Value = 2499.923
$ws.Range("I132").Value = 1761.238
$ws.Range("J132").Value = 5602.4
$ws.Range("K132").Value = 5283.714
$ws.Range("L132").Value = 16807.2
$ws.Range("M132").Value = -2753.714
$ws.Range("N132").Value = -21867.2
$ws.Range("H136").Value = 1404.5122
$ws.Range("I136").Value = 703.7727
$ws.Range("J136").Value = 2215.8948
$ws.Range("K136").Value = 2111.3181
$ws.Range("L136").Value = 6647.6844
$ws.Range("M136").Value = 438.6819
$ws.Range("N136").Value = -11747.6844

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1019.4074
$ws.Range("I5").Value = 548.75
$ws.Range("K5").Value = 1646.25
$ws.Range("M5").Value = -1534.25
$ws.Range("H135").Value = 1019.4074
$ws.Range("I135").Value = 548.75
$ws.Range("K135").Value = 4938.75
$ws.Range("M135").Value = -2403.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1965.0426
$ws.Range("I122").Value = 1672.8235
$ws.Range("J122").Value = 2729.3076
$ws.Range("K122").Value = 5018.470499999999
$ws.Range("L122").Value = 8187.9228
$ws.Range("M122").Value = -2568.470499999999
$ws.Range("N122").Value = -13087.9228
$ws.Range("H132").Value = 2760.5
$ws.Range("I132").Value = 2004.7778
$ws.Range("J132").Value = 4120.8
$ws.Range("K132").Value = 6014.3334
$ws.Range("L132").Value = 12362.4
$ws.Range("M132").Value = -3484.3334
$ws.Range("N132").Value = -17422.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1911.4445
$ws.Range("I22").Value = 1477.7778
$ws.Range("J22").Value = 2345.111
$ws.Range("K22").Value = 1477.7778
$ws.Range("L22").Value = 2345.111
$ws.Range("M22").Value = -1182.7778
$ws.Range("N22").Value = -2935.111
$ws.Range("H27").Value = 1911.4445
$ws.Range("I27").Value = 1477.7778
$ws.Range("J27").Value = 2345.111
$ws.Range("K27").Value = 1477.7778
$ws.Range("L27").Value = 2345.111
$ws.Range("M27").Value = -1370.7778
$ws.Range("N27").Value = -2559.111
$ws.Range("H132").Value = 8730.104499999999
$ws.Range("I132").Value = 10792.533
$ws.Range("J132").Value = 5292.722
$ws.Range("K132").Value = 32377.599
$ws.Range("L132").Value = 15878.166
$ws.Range("M132").Value = -29847.599
$ws.Range("N132").Value = -20938.166
$ws.Range("H136").Value = 25645758
$ws.Range("I136").Value = 4167.1113
$ws.Range("J136").Value = 83339336
$ws.Range("K136").Value = 12501.3339
$ws.Range("L136").Value = 250018008
$ws.Range("M136").Value = -9951.333899999998
$ws.Range("N136").Value = -250023108

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 990.8461
$ws.Range("I81").Value = 896.5
$ws.Range("J81").Value = 1071.7142
$ws.Range("K81").Value = 1793
$ws.Range("L81").Value = 2143.4284
$ws.Range("M81").Value = -732
$ws.Range("N81").Value = -4265.4284
$ws.Range("H84").Value = 990.8461
$ws.Range("I84").Value = 896.5
$ws.Range("J84").Value = 1071.7142
$ws.Range("K84").Value = 8965
$ws.Range("L84").Value = 10717.142
$ws.Range("M84").Value = -3661
$ws.Range("N84").Value = -21325.142
$ws.Range("H132").Value = 1925.6177
$ws.Range("I132").Value = 1514.72
$ws.Range("J132").Value = 3067
$ws.Range("K132").Value = 4544.16
$ws.Range("L132").Value = 9201
$ws.Range("M132").Value = -2014.16
$ws.Range("N132").Value = -14261
$ws.Range("H136").Value = 25901.305
$ws.Range("I136").Value = 39107.348
$ws.Range("J136").Value = 8733.450000000001
$ws.Range("K136").Value = 117322.044
$ws.Range("L136").Value = 26200.35
$ws.Range("M136").Value = -114772.044
$ws.Range("N136").Value = -31300.35
